# MemberData.xlsx test fixture: give the five sample members distinct
# names/emails/postcodes instead of all being copies of "Bob Smith", and
# fix the "grant" -> "Grant" capitalisation typo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (ID 2) - Jill Grant
$ws.Range("D3").Value = "Grant"
$ws.Range("F3").Value = "jill@give.com"
$ws.Range("K3").Value = "W7 5EW"

# Row 4 (ID 3) - John Smith
$ws.Range("C4").Value = "John"
$ws.Range("E4").Value = "John"
$ws.Range("F4").Value = "john@give.com"
$ws.Range("K4").Value = "W1 3QP"

# Row 5 (ID 4) - Cathy Holmes
$ws.Range("B5").Value = "Miss"
$ws.Range("C5").Value = "Cathy"
$ws.Range("D5").Value = "Holmes"
$ws.Range("E5").Value = "Cathy"
$ws.Range("F5").Value = "cathy@give.com"
$ws.Range("K5").Value = "W1 4QP"

# Row 6 (ID 5) - Brian Monroe
$ws.Range("C6").Value = "Brian"
$ws.Range("D6").Value = "Monroe"
$ws.Range("E6").Value = "Brian"
$ws.Range("F6").Value = "brian@give.com"
$ws.Range("K6").Value = "W1 5QP"

# Update the saved selection / active cell shown when the workbook is reopened
$ws.Range("B2:K6").Select()
